$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "BOUNDARY" attack block: columns AK:AR (8 epsilon columns) ---

# 1) Merge the header cell range first (matches how the other attack-group
#    headers, e.g. AC1:AJ1, are merged), then copy formatting from the
#    neighboring "FGSM" header/epsilon cells so the new cells pick up the
#    same bold/border/center style (style index 1) rather than a fresh style.
$ws.Range("AK1:AR1").Merge()
$ws.Range("AJ1").Copy()
$ws.Range("AK1:AR1").PasteSpecial(-4122)
$ws.Range("AK1").Value = "BOUNDARY"

# epsilon row: stored as text labels ("0.01" .. "0.20"), same as the other
# attack blocks - use a leading apostrophe so Excel keeps them as text
# (otherwise "0.10"/"0.20" would lose their trailing zero as a number).
# Set the text values FIRST, then paste the formatting on top, so the
# quote-prefix text entry doesn't leave behind an extra "quotePrefix" style.
$ws.Range("AK2").Value = "'0.01"
$ws.Range("AL2").Value = "'0.02"
$ws.Range("AM2").Value = "'0.03"
$ws.Range("AN2").Value = "'0.04"
$ws.Range("AO2").Value = "'0.05"
$ws.Range("AP2").Value = "'0.07"
$ws.Range("AQ2").Value = "'0.10"
$ws.Range("AR2").Value = "'0.20"

$ws.Range("AJ2").Copy()
$ws.Range("AK2:AR2").PasteSpecial(-4122)

# 2) Numeric metric data (MAE / RMSE / SIM) for each of the three models.

$ws.Range("AK4").Value = 434.262713578542
$ws.Range("AL4").Value = 438.4168499056498
$ws.Range("AM4").Value = 444.4971479860941
$ws.Range("AN4").Value = 450.5092791493734
$ws.Range("AO4").Value = 448.7619357617696
$ws.Range("AP4").Value = 474.8387590408325
$ws.Range("AQ4").Value = 516.8151203282674
$ws.Range("AR4").Value = 690.0286581484477
$ws.Range("AK5").Value = 538.3501106486596
$ws.Range("AL5").Value = 541.4316735494499
$ws.Range("AM5").Value = 552.6875215550439
$ws.Range("AN5").Value = 547.6097236656603
$ws.Range("AO5").Value = 551.7773384544651
$ws.Range("AP5").Value = 597.282424038361
$ws.Range("AQ5").Value = 660.4067212343318
$ws.Range("AR5").Value = 863.6221486541458
$ws.Range("AK6").Value = 0.9991555891704178
$ws.Range("AL6").Value = 0.999150194183373
$ws.Range("AM6").Value = 0.9991195457902678
$ws.Range("AN6").Value = 0.9991261566000184
$ws.Range("AO6").Value = 0.9991158067249729
$ws.Range("AP6").Value = 0.9988355275105271
$ws.Range("AQ6").Value = 0.9986635665302224
$ws.Range("AR6").Value = 0.9976064753285732
$ws.Range("AK7").Value = 368.5963822110494
$ws.Range("AL7").Value = 370.2151553789774
$ws.Range("AM7").Value = 375.5289428393046
$ws.Range("AN7").Value = 386.3836320177714
$ws.Range("AO7").Value = 390.1213578478495
$ws.Range("AP7").Value = 426.0095285987854
$ws.Range("AQ7").Value = 444.4051846758525
$ws.Range("AR7").Value = 704.1476362291971
$ws.Range("AK8").Value = 478.1675443802905
$ws.Range("AL8").Value = 479.8270640448296
$ws.Range("AM8").Value = 483.8208102848614
$ws.Range("AN8").Value = 501.4164901769042
$ws.Range("AO8").Value = 496.1767225970677
$ws.Range("AP8").Value = 539.1085291911966
$ws.Range("AQ8").Value = 566.4028017496709
$ws.Range("AR8").Value = 914.6740281328069
$ws.Range("AK9").Value = 0.9991445502045598
$ws.Range("AL9").Value = 0.9991429928549087
$ws.Range("AM9").Value = 0.9991240140156837
$ws.Range("AN9").Value = 0.9990533151054453
$ws.Range("AO9").Value = 0.9990715664975267
$ws.Range("AP9").Value = 0.9989491958626312
$ws.Range("AQ9").Value = 0.9987819370625872
$ws.Range("AR9").Value = 0.9968723001433576
$ws.Range("AK10").Value = 282.5479021072388
$ws.Range("AL10").Value = 287.4827191543579
$ws.Range("AM10").Value = 304.9597949473063
$ws.Range("AN10").Value = 307.9579650052389
$ws.Range("AO10").Value = 334.3933588091533
$ws.Range("AP10").Value = 345.3904226430257
$ws.Range("AQ10").Value = 426.3090173657735
$ws.Range("AR10").Value = 691.0386287053426
$ws.Range("AK11").Value = 406.6227681327846
$ws.Range("AL11").Value = 412.2370333815308
$ws.Range("AM11").Value = 416.8491417397374
$ws.Range("AN11").Value = 435.8916817664871
$ws.Range("AO11").Value = 453.0843985734363
$ws.Range("AP11").Value = 466.5696390869754
$ws.Range("AQ11").Value = 545.1029101207242
$ws.Range("AR11").Value = 825.7780046604233
$ws.Range("AK12").Value = 0.9993621256692868
$ws.Range("AL12").Value = 0.9993426057079848
$ws.Range("AM12").Value = 0.9993273844989845
$ws.Range("AN12").Value = 0.9992658773192096
$ws.Range("AO12").Value = 0.9992130206408877
$ws.Range("AP12").Value = 0.9991752767832138
$ws.Range("AQ12").Value = 0.9988809465415126
$ws.Range("AR12").Value = 0.9974483385553613
